$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells -------------------------------------------------
# Row 48: I48 0.25 -> 2.5, J48 0.2 -> 2.5
$ws.Range("I48").Value = 2.5
$ws.Range("J48").Value = 2.5

# Row 51: C51 65 -> 55
$ws.Range("C51").Value = 55

# --- Append new rows 61-64 --------------------------------------------------

# Row 61 uses the same formatting pattern as row 54 (full row copy)
$ws.Range("A54:W54").Copy()
$ws.Range("A61:W61").PasteSpecial(-4122)

# Row 62 / Row 63 also use row 54's formatting pattern, except column E,
# which should be blank/empty like E58
$ws.Range("A54:W54").Copy()
$ws.Range("A62:W62").PasteSpecial(-4122)
$ws.Range("A54:W54").Copy()
$ws.Range("A63:W63").PasteSpecial(-4122)

$ws.Range("E58").Copy()
$ws.Range("E62").PasteSpecial(-4122)
$ws.Range("E58").Copy()
$ws.Range("E63").PasteSpecial(-4122)

# Row 64 uses row 6's formatting pattern, except columns A, B, W which follow
# row 60's formatting (index/name/code columns)
$ws.Range("A6:W6").Copy()
$ws.Range("A64:W64").PasteSpecial(-4122)
$ws.Range("A60").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("B60").Copy()
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("W60").Copy()
$ws.Range("W64").PasteSpecial(-4122)

# --- Row 61 values -----------------------------------------------------
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = "Творожный с травами ""Вкусвилл"", 65%, 0,14 кг, пл/с"
$ws.Range("C61").Value = 65
$ws.Range("D61").Value = "Да"
$ws.Range("E61").Value = "Травы"
$ws.Range("F61").Value = "Творожный"
$ws.Range("G61").Value = "Маскарпоне"
$ws.Range("H61").Value = "ВкусВилл"
$ws.Range("I61").Value = 0.14
$ws.Range("J61").Value = 0.14
$ws.Range("K61").Value = 6
$ws.Range("L61").Value = 50
$ws.Range("M61").Value = 0
$ws.Range("N61").Value = 30
$ws.Range("O61").Value = 60
$ws.Range("P61").Value = 0
$ws.Range("Q61").Value = 30
$ws.Range("R61").Value = 25
$ws.Range("S61").Value = 800
$ws.Range("T61").Value = 1000
$ws.Range("U61").Value = 0.7
$ws.Range("V61").Value = 0
$ws.Range("W61").Value = "00-00011195"

# --- Row 62 values -----------------------------------------------------
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = "Творожный ""Вкусвилл"", 65%, 0,25 кг, пл/с"
$ws.Range("C62").Value = 65
$ws.Range("D62").Value = "Да"
$ws.Range("F62").Value = "Творожный"
$ws.Range("G62").Value = "Маскарпоне"
$ws.Range("H62").Value = "ВкусВилл"
$ws.Range("I62").Value = 0.25
$ws.Range("J62").Value = 0.2
$ws.Range("K62").Value = 6
$ws.Range("L62").Value = 50
$ws.Range("M62").Value = 0
$ws.Range("N62").Value = 30
$ws.Range("O62").Value = 60
$ws.Range("P62").Value = 0
$ws.Range("Q62").Value = 30
$ws.Range("R62").Value = 0
$ws.Range("S62").Value = 800
$ws.Range("T62").Value = 1000
$ws.Range("U62").Value = 0.7
$ws.Range("V62").Value = 0
$ws.Range("W62").Value = "00-00011196"

# --- Row 63 values -----------------------------------------------------
$ws.Range("A63").Value = 61
$ws.Range("B63").Value = "Творожный ""Вкусвилл"", 65%, 0,14 кг, пл/с"
$ws.Range("C63").Value = 65
$ws.Range("D63").Value = "Да"
$ws.Range("F63").Value = "Творожный"
$ws.Range("G63").Value = "Маскарпоне"
$ws.Range("H63").Value = "ВкусВилл"
$ws.Range("I63").Value = 0.14
$ws.Range("J63").Value = 0.14
$ws.Range("K63").Value = 6
$ws.Range("L63").Value = 50
$ws.Range("M63").Value = 0
$ws.Range("N63").Value = 30
$ws.Range("O63").Value = 60
$ws.Range("P63").Value = 0
$ws.Range("Q63").Value = 30
$ws.Range("R63").Value = 0
$ws.Range("S63").Value = 800
$ws.Range("T63").Value = 1000
$ws.Range("U63").Value = 0.7
$ws.Range("V63").Value = 0
$ws.Range("W63").Value = "00-00011197"

# --- Row 64 values -----------------------------------------------------
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = "Маскарпоне с шоколадом ""Бонджорно"", 50%, 0,18 кг, пл/с"
$ws.Range("C64").Value = 50
$ws.Range("D64").Value = "Да"
$ws.Range("E64").Value = "Шоколад"
$ws.Range("F64").Value = "Маскарпоне"
$ws.Range("G64").Value = "Маскарпоне"
$ws.Range("H64").Value = "Бонджорно"
$ws.Range("I64").Value = 0.18
$ws.Range("J64").Value = 0.2
$ws.Range("K64").Value = 6
$ws.Range("L64").Value = 50
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = 25
$ws.Range("O64").Value = 50
$ws.Range("P64").Value = 10
$ws.Range("Q64").Value = 0
$ws.Range("R64").Value = 25
$ws.Range("S64").Value = 750
$ws.Range("T64").Value = 1000
$ws.Range("U64").Value = 0.9
$ws.Range("V64").Value = 0
$ws.Range("W64").Value = "00-00011040"

# --- Sheet view bookkeeping --------------------------------------------
$ws.Range("E70").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
